$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.868.03'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.92%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.624.19'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.15%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.42%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.99'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.76%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.48%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.43%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '29.82'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +11.15%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +3.25%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0611'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.51%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.58%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.854.74'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.12%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.625.78'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +3.36%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +6.08%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.91'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +5.02%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.917.38'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.07%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '8.82'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +16.19%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '64.64'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.91%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '244.33'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.56%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0706'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.86%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.29%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +3.60%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.62'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +4.50%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +2.45%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '156.73'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.51%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.61%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.66%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.60'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +3.16%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.995'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.42%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0489'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +3.38%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +5.33%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +3.56%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +3.78%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.425.25'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.21%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.64'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +7.05%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.19%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.04%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.68%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.90%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.556'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +3.42%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0509'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.88%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.00'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.29%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.835'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +4.78%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '54.31'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.97%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '69.23'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +5.05%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +18.32%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.994'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.45%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.43'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.82%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.764.39'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.12%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '88.15'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.83%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +4.19%  '
